$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new incident data (Thach That station)
$ws.Cells.Item(2, 6).Value = "UL_TTT093M_HNI"
$ws.Cells.Item(2, 7).Value = "THACH-HOA-TTT_HNI"
$ws.Cells.Item(2, 10).Value = "POWER_AC_EAS"
$ws.Cells.Item(2, 12).Value = "07/05/2025 13:18:36"
$ws.Cells.Item(2, 20).Value = "Thạch Thất"
$ws.Cells.Item(2, 22).Value = ""
$ws.Cells.Item(2, 27).Value = "Trạm viễn thông loại 3"

# Update row 3 with new incident data (Dan Phuong station)
$ws.Cells.Item(3, 6).Value = "3G_DPG045S_HNI"
$ws.Cells.Item(3, 7).Value = "TRUNG-CHAU-VAN-MON2-11-SMC-DPG_HNI"
$ws.Cells.Item(3, 10).Value = "SITE_OOS"
$ws.Cells.Item(3, 12).Value = "06/05/2025 23:18:18"
$ws.Cells.Item(3, 20).Value = "Đan Phượng"
$ws.Cells.Item(3, 22).Value = "Trạm smc mất điện - 1 - sonnn - 06/05/2025 23:24:34"
$ws.Cells.Item(3, 27).Value = "Trạm viễn thông loại 3"

# Remove the now-obsolete rows 4-8
$ws.Rows("4:8").Delete()

# Adjust column widths to match new content
# (Input values are calibrated so the resulting stored OOXML column width
# is the closest achievable approximation to the target widths of
# 36.7109375 / 14.7109375 / 53.7109375 characters, given this runtime's
# pixel-based width quantization.)
$ws.Columns("G").ColumnWidth = 35.75
$ws.Columns("J").ColumnWidth = 13.75
$ws.Columns("V").ColumnWidth = 52.75
